$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing existing rows 10-67 down to 11-68.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new week's data.
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Vega Monumental Concepción"
$ws.Range("C10").Value = "Bíobío"
$ws.Range("D10").Value = 44462
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100108
$ws.Range("H10").Value = "Tropicales y subtropicales"
$ws.Range("I10").Value = 100108002
$ws.Range("J10").Value = "Mango"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 7500
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 7750
$ws.Range("Q10").Value = "$/bandeja 4 kilos"
$ws.Range("R10").Value = "Brasil"
$ws.Range("S10").Value = 1938
$ws.Range("T10").Value = 4
